# PNAD 2009 - correção nos dados e inicio da analise PNAD 2009
#
# The sheet had two placeholder rows ("situação do domicílio" at row 5 and
# "grandes regiões e unidades da federação" at row 8) whose data values were
# actually recorded one row below where they belonged. The fix removes those
# two empty placeholder rows (shifting every later row up by one, which also
# drops the two now-unused trailing blank rows), and also simplifies the
# header row so the "unnamed: *" placeholder labels become "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header row 2: collapse the "unnamed: 1_level_1" / "unnamed: 5_level_1"
# placeholder labels down to "total".
$ws.Range("B2").Value = "total"
$ws.Range("C2").Value = "total"
$ws.Range("F2").Value = "total"

# Remove the two stray placeholder rows that had no B:F data (the row with
# "situação do domicílio" and, after it shifts up, the row with
# "grandes regiões e unidades da federação"). Deleting them shifts all the
# rows below up into their correct place.
$ws.Rows(5).Delete()
$ws.Rows(7).Delete()
